# Update the daily COVID-19 Valais figures for rows 616-625 (sheet "Feuil1").
# Columns: A=Date, B=Cumul cas positifs (formula), C=Nb nouveaux cas positifs,
#          D=Nb nouvelles admissions, E=Nb nouveaux deces COVID-19,
#          F=Patients intubes, G=Patients hospitalises hors SI,
#          H=Total hospitalisations (formula), I=Nb nouvelles sorties,
#          J=Cumul deces (formula), K=Nb nouveaux deces (formula),
#          L=deces a l'hopital, M=deces extra-hospitaliers.
#
# Only the raw input cells (C, E, F, G) are written here; B, H, J and K are
# driven by the existing shared formulas already in the sheet and recompute
# automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 616: only "Patients hospitalises hors SI" changes.
$ws.Cells.Item(616, 7).Value = 13

# Row 617
$ws.Cells.Item(617, 3).Value = 54
$ws.Cells.Item(617, 7).Value = 12

# Row 618
$ws.Cells.Item(618, 3).Value = 74
$ws.Cells.Item(618, 7).Value = 12

# Row 619
$ws.Cells.Item(619, 3).Value = 61
$ws.Cells.Item(619, 7).Value = 9

# Row 620
$ws.Cells.Item(620, 3).Value = 70
$ws.Cells.Item(620, 7).Value = 9

# Row 621
$ws.Cells.Item(621, 3).Value = 43
$ws.Cells.Item(621, 7).Value = 11

# Row 622
$ws.Cells.Item(622, 3).Value = 103
$ws.Cells.Item(622, 7).Value = 12

# Row 623 - newly filled in (was blank)
$ws.Cells.Item(623, 3).Value = 63
$ws.Cells.Item(623, 5).Value = 3
$ws.Cells.Item(623, 6).Value = 3
$ws.Cells.Item(623, 7).Value = 14

# Row 624 - newly filled in (was blank)
$ws.Cells.Item(624, 3).Value = 80
$ws.Cells.Item(624, 5).Value = 3
$ws.Cells.Item(624, 6).Value = 3
$ws.Cells.Item(624, 7).Value = 14

# Row 625 - newly filled in (was blank)
$ws.Cells.Item(625, 3).Value = 6
$ws.Cells.Item(625, 5).Value = 4
$ws.Cells.Item(625, 6).Value = 3
$ws.Cells.Item(625, 7).Value = 14

$excel.Calculate()
